$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3005427420139313
$ws.Range("B1").Value = 2.519969940185547
$ws.Range("C1").Value = 8.512056350708008
$ws.Range("D1").Value = 2.002977132797241
$ws.Range("E1").Value = 1.16081964969635
